# chore: update Sheets via scheduled runner
# Refreshes cached market-price figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leves across the crafter sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5584.2173
$ws.Range("I40").Value = 4828.5713
$ws.Range("J40").Value = 6759.6665
$ws.Range("K40").Value = 4828.5713
$ws.Range("L40").Value = 6759.6665
$ws.Range("M40").Value = -4653.5713
$ws.Range("N40").Value = -7109.6665
$ws.Range("H80").Value = 11442.5
$ws.Range("I80").Value = 5907.222
$ws.Range("J80").Value = 16977.777
$ws.Range("K80").Value = 17721.666
$ws.Range("L80").Value = 50933.33099999999
$ws.Range("M80").Value = -16723.666
$ws.Range("N80").Value = -52929.33099999999
$ws.Range("H83").Value = 11442.5
$ws.Range("I83").Value = 5907.222
$ws.Range("J83").Value = 16977.777
$ws.Range("K83").Value = 53164.998
$ws.Range("L83").Value = 152799.993
$ws.Range("M83").Value = -48172.998
$ws.Range("N83").Value = -162783.993
$ws.Range("H116").Value = 43762.81
$ws.Range("I116").Value = 73396.92999999999
$ws.Range("J116").Value = 3352.6365
$ws.Range("K116").Value = 73396.92999999999
$ws.Range("L116").Value = 3352.6365
$ws.Range("M116").Value = -69954.92999999999
$ws.Range("N116").Value = -10236.6365
$ws.Range("H132").Value = 2973.5405
$ws.Range("I132").Value = 727.7273
$ws.Range("J132").Value = 21501.5
$ws.Range("K132").Value = 2183.1819
$ws.Range("L132").Value = 64504.5
$ws.Range("M132").Value = 346.8181
$ws.Range("N132").Value = -69564.5
$ws.Range("H135").Value = 1145.375
$ws.Range("I135").Value = 331.75
$ws.Range("J135").Value = 1959
$ws.Range("K135").Value = 2985.75
$ws.Range("L135").Value = 17631
$ws.Range("M135").Value = -450.75
$ws.Range("N135").Value = -22701
$ws.Range("H137").Value = 2838722
$ws.Range("I137").Value = 6076015
$ws.Range("K137").Value = 18228045
$ws.Range("M137").Value = -18225495
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 687.9375
$ws.Range("I45").Value = 652.7
$ws.Range("K45").Value = 652.7
$ws.Range("M45").Value = -275.7
$ws.Range("H61").Value = 575498.5600000001
$ws.Range("I61").Value = 418385.75
$ws.Range("K61").Value = 418385.75
$ws.Range("M61").Value = -418173.75
$ws.Range("H124").Value = 19500
$ws.Range("J124").Value = 19500
$ws.Range("L124").Value = 19500
$ws.Range("N124").Value = -29320
$ws.Range("H136").Value = 575498.5600000001
$ws.Range("I136").Value = 418385.75
$ws.Range("K136").Value = 1255157.25
$ws.Range("M136").Value = -1252607.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2204.1614
$ws.Range("I99").Value = 2495.1538
$ws.Range("J99").Value = 691
$ws.Range("K99").Value = 2495.1538
$ws.Range("L99").Value = 691
$ws.Range("M99").Value = -997.1538
$ws.Range("N99").Value = -3687
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 38.625
$ws.Range("I7").Value = 38.166668
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 38.166668
$ws.Range("L7").Value = 40
$ws.Range("M7").Value = 74.833332
$ws.Range("N7").Value = -266
$ws.Range("H16").Value = 887.5714
$ws.Range("I16").Value = 886.6
$ws.Range("J16").Value = 890
$ws.Range("K16").Value = 886.6
$ws.Range("L16").Value = 890
$ws.Range("M16").Value = -599.6
$ws.Range("N16").Value = -1464
$ws.Range("H31").Value = 1247970.9
$ws.Range("I31").Value = 855.7646999999999
$ws.Range("J31").Value = 2335199.5
$ws.Range("K31").Value = 855.7646999999999
$ws.Range("L31").Value = 2335199.5
$ws.Range("M31").Value = -560.7646999999999
$ws.Range("N31").Value = -2335789.5
$ws.Range("H34").Value = 1247970.9
$ws.Range("I34").Value = 855.7646999999999
$ws.Range("J34").Value = 2335199.5
$ws.Range("K34").Value = 855.7646999999999
$ws.Range("L34").Value = 2335199.5
$ws.Range("M34").Value = -653.7646999999999
$ws.Range("N34").Value = -2335603.5
$ws.Range("H58").Value = 3355.6226
$ws.Range("I58").Value = 4355.6924
$ws.Range("J58").Value = 2392.5925
$ws.Range("K58").Value = 4355.6924
$ws.Range("L58").Value = 2392.5925
$ws.Range("M58").Value = -4152.6924
$ws.Range("N58").Value = -2798.5925
$ws.Range("H99").Value = 27003.691
$ws.Range("I99").Value = 32417
$ws.Range("K99").Value = 32417
$ws.Range("M99").Value = -30919
$ws.Range("H107").Value = 450.0476
$ws.Range("J107").Value = 424.8
$ws.Range("L107").Value = 424.8
$ws.Range("N107").Value = -4264.8
$ws.Range("H113").Value = 887.5714
$ws.Range("I113").Value = 886.6
$ws.Range("J113").Value = 890
$ws.Range("K113").Value = 886.6
$ws.Range("L113").Value = 890
$ws.Range("M113").Value = 1283.4
$ws.Range("N113").Value = -5230
$ws.Range("H126").Value = 27003.691
$ws.Range("I126").Value = 32417
$ws.Range("K126").Value = 97251
$ws.Range("M126").Value = -94781
$ws.Range("H134").Value = 17242824
$ws.Range("I134").Value = 21739994
$ws.Range("J134").Value = 3674
$ws.Range("K134").Value = 65219982
$ws.Range("L134").Value = 11022
$ws.Range("M134").Value = -65217447
$ws.Range("N134").Value = -16092
$ws.Range("H136").Value = 3355.6226
$ws.Range("I136").Value = 4355.6924
$ws.Range("J136").Value = 2392.5925
$ws.Range("K136").Value = 13067.0772
$ws.Range("L136").Value = 7177.7775
$ws.Range("M136").Value = -10517.0772
$ws.Range("N136").Value = -12277.7775
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 763.175
$ws.Range("I5").Value = 493.32
$ws.Range("J5").Value = 1212.9333
$ws.Range("K5").Value = 1479.96
$ws.Range("L5").Value = 3638.7999
$ws.Range("M5").Value = -1367.96
$ws.Range("N5").Value = -3862.7999
$ws.Range("H12").Value = 35.054054
$ws.Range("I12").Value = 22.6
$ws.Range("K12").Value = 67.80000000000001
$ws.Range("M12").Value = 105.2
$ws.Range("H69").Value = 1893
$ws.Range("J69").Value = 10000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622
$ws.Range("H72").Value = 1893
$ws.Range("J72").Value = 10000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112
$ws.Range("H131").Value = 1342.1082
$ws.Range("I131").Value = 1075.8334
$ws.Range("J131").Value = 1469.92
$ws.Range("K131").Value = 3227.5002
$ws.Range("L131").Value = 4409.76
$ws.Range("M131").Value = 1812.4998
$ws.Range("N131").Value = -14489.76
$ws.Range("H135").Value = 763.175
$ws.Range("I135").Value = 493.32
$ws.Range("J135").Value = 1212.9333
$ws.Range("K135").Value = 4439.88
$ws.Range("L135").Value = 10916.3997
$ws.Range("M135").Value = -1904.88
$ws.Range("N135").Value = -15986.3997
$ws.Range("H140").Value = 9350
$ws.Range("I140").Value = 918.64703
$ws.Range("J140").Value = 27266.625
$ws.Range("K140").Value = 2755.94109
$ws.Range("L140").Value = 81799.875
$ws.Range("M140").Value = 2424.05891
$ws.Range("N140").Value = -92159.875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2524.524
$ws.Range("I122").Value = 2934.5
$ws.Range("J122").Value = 2360.5334
$ws.Range("K122").Value = 8803.5
$ws.Range("L122").Value = 7081.600199999999
$ws.Range("M122").Value = -6353.5
$ws.Range("N122").Value = -11981.6002
$ws.Range("H126").Value = 2775.0527
$ws.Range("I126").Value = 2395.875
$ws.Range("J126").Value = 3050.818
$ws.Range("K126").Value = 7187.625
$ws.Range("L126").Value = 9152.454000000002
$ws.Range("M126").Value = -4717.625
$ws.Range("N126").Value = -14092.454
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1283.7778
$ws.Range("I126").Value = 1260.5333
$ws.Range("K126").Value = 3781.5999
$ws.Range("M126").Value = -1311.5999
$ws.Range("H132").Value = 2381.0264
$ws.Range("I132").Value = 1807.5385
$ws.Range("J132").Value = 3623.5833
$ws.Range("K132").Value = 5422.6155
$ws.Range("L132").Value = 10870.7499
$ws.Range("M132").Value = -2892.6155
$ws.Range("N132").Value = -15930.7499
